$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 556.4167
$ws.Range("I92").Value = 574.1111
$ws.Range("J92").Value = 503.33334
$ws.Range("K92").Value = 574.1111
$ws.Range("L92").Value = 503.33334
$ws.Range("M92").Value = 673.8889
$ws.Range("N92").Value = -2999.33334

$ws.Range("H101").Value = 778.6667
$ws.Range("I101").Value = 778.6667
$ws.Range("J101").Value = 0
$ws.Range("K101").Value = 2336.0001
$ws.Range("L101").Value = 0
$ws.Range("M101").ClearContents()
$ws.Range("N101").Value = -714.0001000000002

$ws.Range("H109").Value = 40651.5
$ws.Range("I109").Value = 0
$ws.Range("J109").Value = 40651.5
$ws.Range("K109").Value = 0
$ws.Range("L109").Value = 40651.5
$ws.Range("N109").Value = -43425.5

$ws.Range("I135").Value = 764.4286
$ws.Range("J135").Value = 33334358
$ws.Range("K135").Value = 6879.8574
$ws.Range("L135").Value = 300009222
$ws.Range("M135").Value = -4344.8574
$ws.Range("N135").Value = -300014292

$ws.Range("H137").Value = 4036.2979
$ws.Range("I137").Value = 1353.0883
$ws.Range("J137").Value = 11053.923
$ws.Range("K137").Value = 4059.2649
$ws.Range("L137").Value = 33161.769
$ws.Range("M137").Value = -1509.2649
$ws.Range("N137").Value = -38261.769

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H118").Value = 49202
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 49202
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 49202
$ws.Range("N118").Value = -52516

$ws.Range("H123").Value = 38995.75
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 38995.75
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 38995.75
$ws.Range("N123").Value = -48795.75

$ws.Range("H131").Value = 44154.832
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 44154.832
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 44154.832
$ws.Range("N131").Value = -54234.832

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H23").Value = 6500
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 6500
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 6500
$ws.Range("N23").Value = -7066

$ws.Range("H133").Value = 44198.2
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 44198.2
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 44198.2
$ws.Range("N133").Value = -54318.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2618.5454
$ws.Range("I31").Value = 1207.175
$ws.Range("J31").Value = 3575.4067
$ws.Range("K31").Value = 1207.175
$ws.Range("L31").Value = 3575.4067
$ws.Range("M31").Value = -912.175
$ws.Range("N31").Value = -4165.4067

$ws.Range("H34").Value = 2618.5454
$ws.Range("I34").Value = 1207.175
$ws.Range("J34").Value = 3575.4067
$ws.Range("K34").Value = 1207.175
$ws.Range("L34").Value = 3575.4067
$ws.Range("M34").Value = -1005.175
$ws.Range("N34").Value = -3979.4067

$ws.Range("H58").Value = 1396.8
$ws.Range("I58").Value = 932.63635
$ws.Range("J58").Value = 3585
$ws.Range("K58").Value = 932.63635
$ws.Range("L58").Value = 3585
$ws.Range("M58").Value = -729.63635
$ws.Range("N58").Value = -3991

$ws.Range("H94").Value = 1570.4615
$ws.Range("I94").Value = 970.6667
$ws.Range("J94").Value = 1750.4
$ws.Range("K94").Value = 970.6667
$ws.Range("L94").Value = 1750.4
$ws.Range("M94").Value = -519.6667
$ws.Range("N94").Value = -2652.4

$ws.Range("H132").Value = 63307.434
$ws.Range("I132").Value = 2337.4443
$ws.Range("J132").Value = 282799.4
$ws.Range("K132").Value = 7012.3329
$ws.Range("L132").Value = 848398.2000000001
$ws.Range("M132").Value = -4482.3329
$ws.Range("N132").Value = -853458.2000000001

$ws.Range("H134").Value = 1170187.5
$ws.Range("I134").Value = 2102.75
$ws.Range("J134").Value = 1754229.9
$ws.Range("K134").Value = 6308.25
$ws.Range("L134").Value = 5262689.699999999
$ws.Range("M134").Value = -3773.25
$ws.Range("N134").Value = -5267759.699999999

$ws.Range("H136").Value = 1396.8
$ws.Range("I136").Value = 932.63635
$ws.Range("J136").Value = 3585
$ws.Range("K136").Value = 2797.90905
$ws.Range("L136").Value = 10755
$ws.Range("M136").Value = -247.9090500000002
$ws.Range("N136").Value = -15855

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 1134.9
$ws.Range("I7").Value = 1373.625
$ws.Range("J7").Value = 180
$ws.Range("K7").Value = 4120.875
$ws.Range("L7").Value = 540
$ws.Range("M7").Value = -4008.875
$ws.Range("N7").Value = -764

$ws.Range("H11").Value = 2745.9546
$ws.Range("I11").Value = 2836.238
$ws.Range("J11").Value = 850
$ws.Range("K11").Value = 8508.714
$ws.Range("L11").Value = 2550
$ws.Range("M11").Value = -8368.714
$ws.Range("N11").Value = -2830

$ws.Range("H15").Value = 395.16666
$ws.Range("I15").Value = 85.5
$ws.Range("J15").Value = 550
$ws.Range("K15").Value = 256.5
$ws.Range("L15").Value = 1650
$ws.Range("M15").Value = -116.5
$ws.Range("N15").Value = -1930

$ws.Range("H44").Value = 890.3333
$ws.Range("I44").Value = 751.125
$ws.Range("J44").Value = 2004
$ws.Range("K44").Value = 2253.375
$ws.Range("L44").Value = 6012
$ws.Range("M44").Value = -1855.375
$ws.Range("N44").Value = -6808

$ws.Range("H45").Value = 949.5714
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 949.5714
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 2848.7142
$ws.Range("N45").Value = -3912.7142

$ws.Range("H50").Value = 872.6
$ws.Range("I50").Value = 400
$ws.Range("J50").Value = 1187.6666
$ws.Range("K50").Value = 1200
$ws.Range("L50").Value = 3562.9998
$ws.Range("M50").Value = -719
$ws.Range("N50").Value = -4524.9998

$ws.Range("H53").Value = 872.6
$ws.Range("I53").Value = 400
$ws.Range("J53").Value = 1187.6666
$ws.Range("K53").Value = 1200
$ws.Range("L53").Value = 3562.9998
$ws.Range("M53").Value = -719
$ws.Range("N53").Value = -4524.9998

$ws.Range("H62").Value = 1245
$ws.Range("I62").Value = 1000
$ws.Range("J62").Value = 1490
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 4470
$ws.Range("M62").Value = -2314
$ws.Range("N62").Value = -5842

$ws.Range("H64").Value = 2720.9285
$ws.Range("I64").Value = 1800
$ws.Range("J64").Value = 2972.0908
$ws.Range("K64").Value = 5400
$ws.Range("L64").Value = 8916.2724
$ws.Range("M64").Value = -5130
$ws.Range("N64").Value = -9456.2724

$ws.Range("H65").Value = 1245
$ws.Range("I65").Value = 1000
$ws.Range("J65").Value = 1490
$ws.Range("K65").Value = 9000
$ws.Range("L65").Value = 13410
$ws.Range("M65").Value = -5568
$ws.Range("N65").Value = -20274

$ws.Range("H67").Value = 2720.9285
$ws.Range("I67").Value = 1800
$ws.Range("J67").Value = 2972.0908
$ws.Range("K67").Value = 5400
$ws.Range("L67").Value = 8916.2724
$ws.Range("M67").Value = -4464
$ws.Range("N67").Value = -10788.2724

$ws.Range("H74").Value = 14999.75
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 14999.75
$ws.Range("K74").Value = 0
$ws.Range("L74").Value = 44999.25
$ws.Range("N74").Value = -47121.25

$ws.Range("H75").Value = 3176.2856
$ws.Range("I75").Value = 783.6667
$ws.Range("J75").Value = 4133.3335
$ws.Range("K75").Value = 2351.0001
$ws.Range("L75").Value = 12400.0005
$ws.Range("M75").Value = -1353.0001
$ws.Range("N75").Value = -14396.0005

$ws.Range("H77").Value = 14999.75
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 14999.75
$ws.Range("K77").Value = 0
$ws.Range("L77").Value = 134997.75
$ws.Range("N77").Value = -145605.75

$ws.Range("H78").Value = 3176.2856
$ws.Range("I78").Value = 783.6667
$ws.Range("J78").Value = 4133.3335
$ws.Range("K78").Value = 7053.0003
$ws.Range("L78").Value = 37200.0015
$ws.Range("M78").Value = -2061.0003
$ws.Range("N78").Value = -47184.0015

$ws.Range("H131").Value = 3551.476
$ws.Range("I131").Value = 11585.444
$ws.Range("J131").Value = 1360.3939
$ws.Range("K131").Value = 34756.33199999999
$ws.Range("L131").Value = 4081.1817
$ws.Range("M131").Value = -29716.33199999999
$ws.Range("N131").Value = -14161.1817

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H21").Value = 6833.3335
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 6833.3335
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 6833.3335
$ws.Range("N21").Value = -7179.3335

$ws.Range("H30").Value = 6833.3335
$ws.Range("I30").Value = 0
$ws.Range("J30").Value = 6833.3335
$ws.Range("K30").Value = 0
$ws.Range("L30").Value = 6833.3335
$ws.Range("N30").Value = -7043.3335

$ws.Range("H122").Value = 1046.5555
$ws.Range("I122").Value = 1039.8667
$ws.Range("J122").Value = 1080
$ws.Range("K122").Value = 3119.6001
$ws.Range("L122").Value = 3240
$ws.Range("M122").Value = -669.6001000000001
$ws.Range("N122").Value = -8140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 23438
$ws.Range("I2").Value = 1500

$ws.Range("H21").Value = 20000
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 20000
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 20000
$ws.Range("N21").Value = -20348

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 7500
$ws.Range("I20").Value = 5000
$ws.Range("J20").Value = 10000
$ws.Range("K20").Value = 5000
$ws.Range("L20").Value = 10000
$ws.Range("M20").Value = -4760
$ws.Range("N20").Value = -10480

$ws.Range("H138").Value = 45200
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 45200
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 45200
$ws.Range("N138").Value = -55480
